$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 535 (existing rows 535-601 shift down to 537-603)
$ws.Rows("535:536").Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown)

# New row 535 - same as old row 535 except D, M, N, O, P, S
$ws.Range("A535").Value = 2
$ws.Range("B535").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C535").Value = "Coquimbo"
$ws.Range("D535").Value = 44748
$ws.Range("D535").NumberFormat = $ws.Range("D537").NumberFormat
$ws.Range("E535").Value = 4
$ws.Range("F535").Value = "Fruta"
$ws.Range("G535").Value = 100102
$ws.Range("H535").Value = "Cítricos"
$ws.Range("I535").Value = 100102003
$ws.Range("J535").Value = "Limón"
$ws.Range("K535").Value = "Sin especificar"
$ws.Range("L535").Value = "1a amarillo"
$ws.Range("M535").Value = 600
$ws.Range("N535").Value = 2300
$ws.Range("O535").Value = 2500
$ws.Range("P535").Value = 2400
$ws.Range("Q535").Value = "$/malla 16 kilos"
$ws.Range("R535").Value = "Provincia de Limarí"
$ws.Range("S535").Value = 150
$ws.Range("T535").Value = 16

# New row 536 - same as old row 536 except D, M, N, O, P, S
$ws.Range("A536").Value = 2
$ws.Range("B536").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C536").Value = "Coquimbo"
$ws.Range("D536").Value = 44748
$ws.Range("D536").NumberFormat = $ws.Range("D537").NumberFormat
$ws.Range("E536").Value = 4
$ws.Range("F536").Value = "Fruta"
$ws.Range("G536").Value = 100102
$ws.Range("H536").Value = "Cítricos"
$ws.Range("I536").Value = 100102003
$ws.Range("J536").Value = "Limón"
$ws.Range("K536").Value = "Sin especificar"
$ws.Range("L536").Value = "2a amarillo"
$ws.Range("M536").Value = 600
$ws.Range("N536").Value = 1300
$ws.Range("O536").Value = 1500
$ws.Range("P536").Value = 1400
$ws.Range("Q536").Value = "$/malla 16 kilos"
$ws.Range("R536").Value = "Provincia de Limarí"
$ws.Range("S536").Value = 88
$ws.Range("T536").Value = 16
